$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Credits / hours / activation date fields
Replace-Text "Créditos-aula: 4" "Créditos-aula: 2"
Replace-Text "Carga horária: 60 h" "Carga horária: 30 h"
Replace-Text "Ativação: 01/01/2022" "Ativação: 01/01/2025"

# "Programa resumido" section - Portuguese summary paragraph
Replace-Text "Principais fatores condicionantes atuantes em processos como erosão, movimento gravitacional de massa e contaminação de materiais geológicos e água subterrânea. Técnicas de investigação, monitoramento, contenção e recuperação em geotecnia ambiental." "Principais fatores condicionantes atuantes em processos da dinâmica superficial. Técnicas de contenção e recuperação em geotecnia ambiental."

# "Programa resumido" section - English summary paragraph
Replace-Text "Main aspects acting in processes such as erosion, gravitational mass movements and contamination of geological materials and groundwater. Investigation, monitoring, containment and recuperation techniques in environmental geotechnics." "Main aspects acting in processes of the surficial dynamics. Containment and recuperation techniques in environmental geotechnics"

# "Programa" section - Portuguese full paragraph
Replace-Text "Aspectos conceituais e abrangência do tema. Investigação e monitoramento geoambiental. Mapas geoambientais e geotécnicos. Movimentos gravitacionais de massa: conceitos, tipos e mecanismos. Análise de estabilidade de taludes e técnicas de estabilização de encostas naturais. Processos erosivos. Estimativas de perda de solo e técnicas de mitigação de processos erosivos antrópicos. Contaminação de materiais geológicos e água subterrânea. Transporte de contaminantes na zona saturada e não saturada do solo. Barragens de rejeitos e aterros de resíduos. Geossintéticos" "Introdução a mecânica dos solos ambiental (granulometria, limites de consistência dos solos, compactação, percolação e condutividade hidráulica)Movimentos gravitacionais de massa: conceitos, tipos e mecanismos. Análise de estabilidade de taludes e técnicas de estabilização Processos erosivos: conceitos, tipos e mecanismos. Estimativas de perda de solo e técnicas de mitigação de processos erosivos antrópicosColapso e subsidênciaAterros de resíduos (seleção de áreas para aterros sanitários, métodos construtivos, ABNT 8419 e ABNT 15849)Levantamento de estudos e projetos geotécnicos aplicados a Estudos AmbientaisAplicação de geossintéticos em obras ambientaisAulas práticas: parâmetros de mecânica dos solos A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina"

# "Programa" section - English full paragraph
Replace-Text "Conceptual aspects and scope of the theme. Geoenvironmental investigation and monitoring. Geoenvironmental and geotechnical maps. Gravitational mass movements: concepts, types and mechanics. Slope stability analysis and natural slope stabilization techniques. Erosive processes. Estimates of soil loss and mitigation techniques for anthropic erosive processes. Contamination of geological materials and groundwater. Transport of contaminants in saturated and unsaturated zone of the soil. Tailings dams and waste landfills. Geosynthetics." "Introduction to environmental soil mechanics (granulometry, soil consistency limits, compaction, percolation and hydraulic conductivity)Gravitational mass movements: concepts, types and mechanisms. Slope stability analysis and stabilization techniquesErosive processes: concepts, types and mechanisms. Soil loss estimative and mitigation techniques for anthropogenic erosion processesCollapse and subsidenceWaste landfills (selection of areas for landfills, construction methods, ABNT 8419 and ABNT 15849)Survey of geotechnical studies and projects applied to Environmental StudiesApplication of geosynthetics in environmental projectsExperimental classes: soil mechanics parametersThe discipline may have didactic trips to complement the content of the discipline"
